$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns being written so Excel
# does not auto-convert numeric-looking strings (e.g. "1.00", "606.92")
# or padded percentage strings into numbers and lose their exact text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.336.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.570.33'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.88'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.568.75'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.137'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.84'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.175.54'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.74%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.573.27'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.384.40'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.45'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.59%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.86'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '431.61'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.11%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.76'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.715.16'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000119'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.51'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.21%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.17'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.95'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.52'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.56%  '
$ws.Range("B33").Value = 'RenzoRestakedETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.564.35'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.18%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.39%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.73'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.63'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '175.27'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0853'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.11%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.890'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.00'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.35%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.23'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.15'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.48'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.02%  '
